$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: new daily entry (No=6, Date=2022-01-10, App=RPA RLOGIC, Task, 100%, Completed) ---
# Copy formats first so number formats (date / percentage) match the rest of the table.
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("E9").Copy()
$ws.Range("E10:E12").PasteSpecial(-4122)

$ws.Range("A10").Value = 6
$ws.Range("B10").Value = 44571
$ws.Range("C10").Value = "RPA RLOGIC"
$ws.Range("D10").Value = "1. Expenses data has been added to the P&L report as separate sheet like callwise data sheet"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = "Completed"

# --- Row 11: second bullet of the same day's task ---
$ws.Range("D11").Value = "2. Deleted the backup files of the Rlogic daily reports from Sep2021 to Dec2021"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = "Completed"

# --- Row 12: third bullet of the same day's task ---
$ws.Range("D12").Value = "3. Deleted the log files from Sep2021 to Dec2021"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = "Completed"

$excel.CutCopyMode = $false

# --- Update the active selection left by the author at D21 ---
$ws.Range("D21").Select() | Out-Null
